$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$f389 = @'
OK
###
Log Config LLM:{'last_update': '2024-11-17 20:22:21', 'llm': 'openai', 'model_llm': 'gpt-4o', 'embbeder': 'openai', 'model_embedder': 'text-embedding-3-large', 'chunk_size': 1000, 'chunk_overlap': 200, 'total_chunks': 1023}
'@

$f393 = @'
Proses penyiapan dokumen berhasil diselesaikan dan embeddings berhasil disimpan pada vector database.
###
llm:openai
###
model_llm:gpt-4o
###
embbeder:openai
###
model_embedder:text-embedding-3-large
###
chunk_size:1000
###
chunk_overlap:200
###
total_chunks:1000
'@

$ws.Cells.Item(386, 1).Value = 'V4S75B'
$ws.Cells.Item(386, 2).Value = '2024-11-22 19:44:17'
$ws.Cells.Item(386, 3).Value = 'GET /'
$ws.Cells.Item(386, 4).Value = 401
$ws.Cells.Item(386, 5).Value = $false
$ws.Cells.Item(386, 6).Value = 'Eitss... mau ngapain? Akses terbatas!'

$ws.Cells.Item(387, 1).Value = '5IKUQO'
$ws.Cells.Item(387, 2).Value = '2024-11-22 19:44:17'
$ws.Cells.Item(387, 3).Value = 'GET /favicon.ico'
$ws.Cells.Item(387, 4).Value = 404
$ws.Cells.Item(387, 5).Value = $false
$ws.Cells.Item(387, 6).Value = 'Not Found'

$ws.Cells.Item(388, 1).Value = 'HDEL6D'
$ws.Cells.Item(388, 2).Value = '2024-11-22 19:44:31'
$ws.Cells.Item(388, 3).Value = 'GET /'
$ws.Cells.Item(388, 4).Value = 200
$ws.Cells.Item(388, 5).Value = $true
$ws.Cells.Item(388, 6).Value = 'API Virtual Assistant Undiksha 2024-11-22 19:44:31'

$ws.Cells.Item(389, 1).Value = 'SWSD4G'
$ws.Cells.Item(389, 2).Value = '2024-11-22 19:45:16'
$ws.Cells.Item(389, 3).Value = 'GET /checkmodel'
$ws.Cells.Item(389, 4).Value = 200
$ws.Cells.Item(389, 5).Value = $true
$ws.Cells.Item(389, 6).Value = $f389

$ws.Cells.Item(390, 1).Value = 'PN9IYD'
$ws.Cells.Item(390, 2).Value = '2024-11-22 19:45:34'
$ws.Cells.Item(390, 3).Value = 'POST /setup'
$ws.Cells.Item(390, 4).Value = 400
$ws.Cells.Item(390, 5).Value = $false
$ws.Cells.Item(390, 6).Value = 'Embedder harus ''openai'' atau ''ollama''.'

$ws.Cells.Item(391, 1).Value = 'GP1KBH'
$ws.Cells.Item(391, 2).Value = '2024-11-22 19:48:45'
$ws.Cells.Item(391, 3).Value = 'POST /setup'
$ws.Cells.Item(391, 4).Value = 400
$ws.Cells.Item(391, 5).Value = $false
$ws.Cells.Item(391, 6).Value = 'Embedder harus ''openai'' atau ''ollama''.'

$ws.Cells.Item(392, 1).Value = 'EG78BY'
$ws.Cells.Item(392, 2).Value = '2024-11-22 19:48:52'
$ws.Cells.Item(392, 3).Value = 'POST /setup'
$ws.Cells.Item(392, 4).Value = 400
$ws.Cells.Item(392, 5).Value = $false
$ws.Cells.Item(392, 6).Value = 'Model Embedder untuk ''openai'' harus salah satu dari [''text-embedding-3-large'', ''text-embedding-3-small''].'

$ws.Cells.Item(393, 1).Value = 'QKF843'
$ws.Cells.Item(393, 2).Value = '2024-11-22 19:49:52'
$ws.Cells.Item(393, 3).Value = 'POST /setup'
$ws.Cells.Item(393, 4).Value = 200
$ws.Cells.Item(393, 5).Value = $true
$ws.Cells.Item(393, 6).Value = $f393

